# Add anonymous transition and deferred events.
#
# The "Defer event and anonymous transition" slide (slide 12) had an
# extra straight-arrow connector (id=61, "直線矢印コネクタ 28") linking
# shapes 44 -> 58 that shouldn't be there once the anonymous-transition
# state (and its own connectors) were added. Remove that stray connector.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(12)

for ($i = $s.Shapes.Count; $i -ge 1; $i--) {
    $shape = $s.Shapes.Item($i)
    if ($shape.Id -eq 61) {
        $shape.Delete()
        break
    }
}
